$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "CARRION LAZARO MICHAEL LUIS"
$ws.Range("B2").Value = 165

# Row 3
$ws.Range("A3").Value = "ARRUNATEGUI ESPINOZA JOVANNY"
$ws.Range("B3").Value = 163

# Row 4
$ws.Range("A4").Value = "PAZ ANASTACIO JUANITA ROSA"
$ws.Range("B4").Value = 146

# Row 5
$ws.Range("A5").Value = "NIMA CARMEN KAREN DEL MILAGRO"
$ws.Range("B5").Value = 146

# Row 6
$ws.Range("A6").Value = "ESPINOZA VALDIVIEZO JUNIOR RICARDO"
$ws.Range("B6").Value = 135

# Row 7
$ws.Range("A7").Value = "ALZAMORA CHERRES SIRLEY YASMIN"
$ws.Range("B7").Value = 130

# Row 8
$ws.Range("A8").Value = "PULACHE LAZO VILMA YOHANA"
$ws.Range("B8").Value = 129

# Row 9
$ws.Range("A9").Value = "DOMINGUEZ CUEVA MERLING DEL JESUS YOLINDA"
$ws.Range("B9").Value = 121

# Row 10
$ws.Range("A10").Value = "LILIAN ROXANA VEGA GARCÍA"
$ws.Range("B10").Value = 120

# Row 11
$ws.Range("A11").Value = "NAVARRO JUAREZ LIDIA"
$ws.Range("B11").Value = 116
